$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.477.95"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").Value = "1.567.42"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.38"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3689"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.25"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3385"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.143"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07503"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.011"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.977"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").Value = "1.569.28"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("E17").Value = "  +0.66%  "

$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06757"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.373"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.24%  "

$ws.Range("D24").Value = "22.461.52"
$ws.Range("E24").Value = "  +0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.396"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.35%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.630"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.88"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.39"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.52%  "

$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.59"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").Value = "1.744.01"
$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.059"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.194"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.84%  "

$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.720"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08313"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02474"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2283"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.339"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06469"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.399"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.23"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6206"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.89"
$ws.Range("D45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.770"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5829"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.17%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.051"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.71"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.229"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07336"
$ws.Range("D51").ClearFormats()
